$wb = $excel.ActiveWorkbook

# --- Sheet1: "LoginData" ---
$ws1 = $wb.Worksheets.Item(1)

# Remove the old row 3 (user2@test / pass456) - data now fits in 2 rows.
$ws1.Rows.Item(3).Delete()

# Replace row 2 login data with the new sample credentials.
$ws1.Range("A2").Value = "testautomationram@gmail.com"
$ws1.Range("B2").Value = "Test@1234"

# Turn the new values into (mailto:) hyperlinks, like Excel's AutoFormat
# does when you type an "@"-containing value into a cell.
$ws1.Hyperlinks.Add($ws1.Range("A2"), "mailto:testautomationram@gmail.com")
$ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:Test@1234")

# Widen column A to fit the longer e-mail address.
$ws1.Columns.Item(1).ColumnWidth = 30.28515625

# Move the sheet's saved selection.
$ws1.Range("A10").Select()

# --- Sheet2: "InvalidLoginData" ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "bad1@gmail.com"
$ws2.Range("A3").Value = "bad2@gmail.com"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:bad1@gmail.com")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:bad2@gmail.com")

# Move the sheet's saved selection.
$ws2.Range("B15").Select()
